$d = $word.ActiveDocument

# The two trailing boilerplate paragraphs ("Ver no Jupiter Salvar em pdf
# Salvar em docx" and the "(c) 2020 ..." copyright line) were removed from
# the bottom of the page, right after "LOB1019: Física II (Requisito
# fraco)". Find the exact run of text, including the two paragraph marks
# that separate/terminate those paragraphs, and delete it so the following
# (still-empty) page-break paragraph simply reattaches to what precedes it.
$range = $d.Content
$found = $range.Find.Execute(
    "Ver no Jupiter Salvar em pdf Salvar em docx^p© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution^p",
    $false, $false, $false, $false, $false, $true, 1, $false, "", 2
)

if ($found) {
    $range.Delete()
}
